$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Age Group
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Age Group")

$ws.Range("B2").Value = 47972
$ws.Range("C2").Value = 4121
$ws.Range("E2").Value = 8.18
$ws.Range("F2").Value = 7.72

$ws.Range("B3").Value = 90439
$ws.Range("C3").Value = 9055
$ws.Range("E3").Value = 15.42
$ws.Range("F3").Value = 16.97

$ws.Range("B4").Value = 91074
$ws.Range("C4").Value = 8491
$ws.Range("E4").Value = 15.53
$ws.Range("F4").Value = 15.91

$ws.Range("B5").Value = 88199
$ws.Range("C5").Value = 8788
$ws.Range("E5").Value = 15.04
$ws.Range("F5").Value = 16.47

$ws.Range("B6").Value = 95335
$ws.Range("C6").Value = 8232
$ws.Range("E6").Value = 16.25
$ws.Range("F6").Value = 15.42
$ws.Range("G6").Value = 5.13

$ws.Range("B7").Value = 83607
$ws.Range("C7").Value = 6311
$ws.Range("D7").Value = 410
$ws.Range("E7").Value = 14.25
$ws.Range("F7").Value = 11.82
$ws.Range("G7").Value = 15.82

$ws.Range("B8").Value = 52487
$ws.Range("C8").Value = 3968
$ws.Range("D8").Value = 636
$ws.Range("E8").Value = 8.949999999999999
$ws.Range("F8").Value = 7.43
$ws.Range("G8").Value = 24.54

$ws.Range("B9").Value = 34196
$ws.Range("C9").Value = 4372
$ws.Range("D9").Value = 1339
$ws.Range("E9").Value = 5.83
$ws.Range("F9").Value = 8.19
$ws.Range("G9").Value = 51.66

$ws.Range("B10").Value = 3280
$ws.Range("E10").Value = 0.5600000000000001

# ---------------------------------------------------------------------------
# Sheet: Gender
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Gender")

$ws.Range("B2").Value = 336036
$ws.Range("C2").Value = 27884
$ws.Range("D2").Value = 1261
$ws.Range("E2").Value = 57.29

$ws.Range("B3").Value = 238840
$ws.Range("C3").Value = 24567
$ws.Range("D3").Value = 1277
$ws.Range("E3").Value = 40.72
$ws.Range("F3").Value = 46.03

$ws.Range("B4").Value = 11713
$ws.Range("C4").Value = 919
$ws.Range("E4").Value = 2

# ---------------------------------------------------------------------------
# Sheet: Race
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Race")

$ws.Range("B2").Value = 4962
$ws.Range("C2").Value = 748
$ws.Range("E2").Value = 0.85
$ws.Range("F2").Value = 1.4

$ws.Range("B3").Value = 51472
$ws.Range("C3").Value = 6311
$ws.Range("D3").Value = 371
$ws.Range("E3").Value = 8.77
$ws.Range("F3").Value = 11.82
$ws.Range("G3").Value = 14.31

$ws.Range("B4").Value = 79965
$ws.Range("C4").Value = 9259
$ws.Range("D4").Value = 344
$ws.Range("E4").Value = 13.63
$ws.Range("F4").Value = 17.35
$ws.Range("G4").Value = 13.27

$ws.Range("B5").Value = 107037
$ws.Range("C5").Value = 12322
$ws.Range("D5").Value = 182
$ws.Range("E5").Value = 18.25
$ws.Range("F5").Value = 23.09
$ws.Range("G5").Value = 7.02

$ws.Range("B6").Value = 343153
$ws.Range("C6").Value = 24730
$ws.Range("D6").Value = 1682
$ws.Range("E6").Value = 58.5
$ws.Range("F6").Value = 46.34
$ws.Range("G6").Value = 64.89

# ---------------------------------------------------------------------------
# Sheet: Ethnicity
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Ethnicity")

$ws.Range("B2").Value = 19569
$ws.Range("C2").Value = 5971
$ws.Range("E2").Value = 3.34
$ws.Range("F2").Value = 11.19
$ws.Range("G2").Value = 1.97

$ws.Range("B3").Value = 225000
$ws.Range("C3").Value = 17907
$ws.Range("D3").Value = 1367
$ws.Range("E3").Value = 38.36
$ws.Range("F3").Value = 33.55
$ws.Range("G3").Value = 52.74

$ws.Range("B4").Value = 342020
$ws.Range("C4").Value = 29492
$ws.Range("D4").Value = 1174
$ws.Range("E4").Value = 58.31
$ws.Range("F4").Value = 55.26
$ws.Range("G4").Value = 45.29

Write-Output "applied indiana covid data updates"
